$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "189"
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "446016.00"

$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = "1002"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3183764.33"

$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value = "414"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1683698.25"

$ws.Range("C5").NumberFormat = "@"
$ws.Range("C5").Value = "117"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "562128.09"

$ws.Range("C9").NumberFormat = "@"
$ws.Range("C9").Value = "57"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "147928.41"

$ws.Range("C10").NumberFormat = "@"
$ws.Range("C10").Value = "357"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1267688.71"

$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = "145"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "618601.77"

$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = "35"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "175120.00"

$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = "16"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "41500.00"

$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = "426"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1343086.23"

$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = "153"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "655702.10"

$ws.Range("C33").NumberFormat = "@"
$ws.Range("C33").Value = "108"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "307173.00"

$ws.Range("C34").NumberFormat = "@"
$ws.Range("C34").Value = "566"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1861379.47"

$ws.Range("C35").NumberFormat = "@"
$ws.Range("C35").Value = "226"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1134788.11"

$ws.Range("C36").NumberFormat = "@"
$ws.Range("C36").Value = "73"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "397894.00"

$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = "101"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "283768.17"

$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = "587"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2031583.52"

$ws.Range("C52").NumberFormat = "@"
$ws.Range("C52").Value = "262"
$ws.Range("D52").NumberFormat = "@"
$ws.Range("D52").Value = "1142878.76"

$ws.Range("C55").NumberFormat = "@"
$ws.Range("C55").Value = "23"
$ws.Range("D55").NumberFormat = "@"
$ws.Range("D55").Value = "68220.65"

$ws.Range("C82").NumberFormat = "@"
$ws.Range("C82").Value = "887"
$ws.Range("D82").NumberFormat = "@"
$ws.Range("D82").Value = "2841512.26"
